$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate column A (rows 1-5) and row 1 (A1:D1) with the same style used
# elsewhere in the sheet (style index 1 -> centered, Times New Roman).
# Copy formatting from an existing style-1 cell (C2) onto the new range,
# then set the one cell that also carries a value (A5 = 5).
$ws.Range("C2").Copy()
$ws.Range("A1:A5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B1:D1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A5").Value = 5

# Move the selection to match the committed sheet view.
[void]$ws.Range("B5").Select()
